$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column D and update each 'Flavors' cell
# value from "['n/a']" to "[]" for every data row (rows 2-7).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "['n/a']") {
        $cell.Value = "[]"
    }
}
